{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst targetParagraph = paragraphs.items[0];\nconst newText = \"1. Les administrations de trois\\n2. villes voisines: A, B et C ont d\u00e9cid\u00e9 \\n3. De construire un a\u00e9roport et de diviser les co\u00fbts de ce\\n4. Projet.\\n5. L\u2019endroit que convient le plus est \\n6.celui dont la somme les distances entre chaque\\n7. Ville et l\u2019a\u00e9roport est la plus petite \\n8. possible. Les experts qui sont responsables \\n9. de ce travail, ont d\u00e9velopp\u00e9 un mod\u00e8le \\n10. Pour avoir une premi\u00e8re id\u00e9e o\u00f9 cet a\u00e9roport\\n11. Pourrait \u00eatre plac\u00e9. Il peuvent utiliser \\n12. des clous und anneau grand et une \\n13. corde longue.\\n14. Expliquez comment les experts peuvent utiliser \\n15. ces mat\u00e9riaux pour trouver une approximation\\n16. de l\u2019endroit optimal pour l\u2019a\u00e9roport. Imaginez\\n17. que les villes sont situ\u00e9es \\n18. aux angles d\u2019un triangle \\n19. qui est \u00e9videmment dessin\u00e9 \u00e0 l\u2019\u00e9chelle \\n20. dans cette figure. Ceci est une possibilit\u00e9: \\n21. La roue commence au premier clou,\\n22.  va dans l\u2019anneau, Puis va autour l\u2019autre clou\\n23. , autour le dernier clou, \\n24.  Et encore dans l\u2019anneau et maintenant on peut tirer \\n25. la corde pour trouver l\u2019endroit \\n26. qu\u2019on cherche. Pour y arriver, \\n27.il faut d\u00e9placer la corde un peu,\\n28.Parce qu\u2019il y a de la r\u00e9sistance \u00e0 cause des\\n29. Mat\u00e9riaux qu\u2019on utilise, mais\\n30. \u00e0 la fin, vous arriviez \u00e0 la position,\\n31. De laquelle l\u2019anneau ne bouge plus,\\n32. qui est plus ou moins celle-ci. \\n33. Vous voyez bien que les angles destrois distances\\n34. Entre les clous et l\u2019anneau \\n35. sont plus au moins \u00e9gaux \u00e0 120 d\u00e9gr\u00e9es\\n36. Ce qui est 1/3 de l\u2019angle plein.\\n37. Ceci est l\u2019endroit qu\u2019on a cherch\u00e9:\\n38. -40. La somme minimale des distances entre les villes et l\u2019a\u00e9roport.\\n41. Musique \\n\";\n\n// Replace the run text in place (preserving paragraph formatting / pPr)\n// while keeping the embedded line breaks as literal newline characters\n// inside a single run, matching the source document's structure.\nconst range = targetParagraph.getRange(\"Content\");\nrange._omSet(\"Text\", newText, \"Range\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$p = $d.Paragraphs(1)\n$newText = @'\n1. Les administrations de trois\n2. villes voisines: A, B et C ont d\u00e9cid\u00e9 \n3. De construire un a\u00e9roport et de diviser les co\u00fbts de ce\n4. Projet.\n5. L\u2019endroit que convient le plus est \n6.celui dont la somme les distances entre chaque\n7. Ville et l\u2019a\u00e9roport est la plus petite \n8. possible. Les experts qui sont responsables \n9. de ce travail, ont d\u00e9velopp\u00e9 un mod\u00e8le \n10. Pour avoir une premi\u00e8re id\u00e9e o\u00f9 cet a\u00e9roport\n11. Pourrait \u00eatre plac\u00e9. Il peuvent utiliser \n12. des clous und anneau grand et une \n13. corde longue.\n14. Expliquez comment les experts peuvent utiliser \n15. ces mat\u00e9riaux pour trouver une approximation\n16. de l\u2019endroit optimal pour l\u2019a\u00e9roport. Imaginez\n17. que les villes sont situ\u00e9es \n18. aux angles d\u2019un triangle \n19. qui est \u00e9videmment dessin\u00e9 \u00e0 l\u2019\u00e9chelle \n20. dans cette figure. Ceci est une possibilit\u00e9: \n21. La roue commence au premier clou,\n22.  va dans l\u2019anneau, Puis va autour l\u2019autre clou\n23. , autour le dernier clou, \n24.  Et encore dans l\u2019anneau et maintenant on peut tirer \n25. la corde pour trouver l\u2019endroit \n26. qu\u2019on cherche. Pour y arriver, \n27.il faut d\u00e9placer la corde un peu,\n28.Parce qu\u2019il y a de la r\u00e9sistance \u00e0 cause des\n29. Mat\u00e9riaux qu\u2019on utilise, mais\n30. \u00e0 la fin, vous arriviez \u00e0 la position,\n31. De laquelle l\u2019anneau ne bouge plus,\n32. qui est plus ou moins celle-ci. \n33. Vous voyez bien que les angles destrois distances\n34. Entre les clous et l\u2019anneau \n35. sont plus au moins \u00e9gaux \u00e0 120 d\u00e9gr\u00e9es\n36. Ce qui est 1/3 de l\u2019angle plein.\n37. Ceci est l\u2019endroit qu\u2019on a cherch\u00e9:\n38. -40. La somme minimale des distances entre les villes et l\u2019a\u00e9roport.\n41. Musique \n'@\n$newText = $newText + \"`n\"\n$p.Range.Text = $newText\n"}
